# Column D prices are plain text in the workbook (coinranking price strings).
# Values that look numeric (single decimal point, e.g. "216.56") must be
# written with a leading apostrophe so Excel stores them as text rather than
# coercing them to a number - matching values with thousands separators
# (e.g. "28.007.88") are already unambiguous text and need no prefix.
$updates = @(
    @{ Row = 2;  D = "28.007.88";    E = "  +3.07%  " },
    @{ Row = 3;  D = "1.686.67";     E = "  +0.20%  " },
    @{ Row = 4;  D = $null;          E = "  -0.16%  " },
    @{ Row = 5;  D = "'216.56";      E = "  +0.48%  " },
    @{ Row = 6;  D = "'0.518";       E = "  +0.02%  " },
    @{ Row = 7;  D = $null;          E = "  -0.18%  " },
    @{ Row = 8;  D = "'24.05";       E = "  +6.25%  " },
    @{ Row = 9;  D = $null;          E = "  +1.77%  " },
    @{ Row = 10; D = "'0.0625";      E = "  +0.21%  " },
    @{ Row = 11; D = $null;          E = "  -0.73%  " },
    @{ Row = 12; D = "1.924.31";     E = "  +0.13%  " },
    @{ Row = 13; D = "1.688.05";     E = "  +0.36%  " },
    @{ Row = 14; D = $null;          E = "  -0.22%  " },
    @{ Row = 15; D = $null;          E = "  +0.09%  " },
    @{ Row = 16; D = "'66.82";       E = "  -0.03%  " },
    @{ Row = 17; D = "'250.59";      E = "  +6.08%  " },
    @{ Row = 18; D = "27.980.16";    E = "  +2.96%  " },
    @{ Row = 19; D = $null;          E = "  +0.24%  " },
    @{ Row = 20; D = $null;          E = "  -3.23%  " },
    @{ Row = 21; D = $null;          E = "  -0.06%  " },
    @{ Row = 22; D = "'4.53";        E = "  -0.47%  " },
    @{ Row = 23; D = $null;          E = "  -0.10%  " },
    @{ Row = 24; D = $null;          E = "  -2.33%  " },
    @{ Row = 25; D = "'147.47";      E = "  +0.40%  " },
    @{ Row = 26; D = "'7.34";        E = "  -0.90%  " },
    @{ Row = 27; D = "'16.48";       E = "  +0.77%  " },
    @{ Row = 28; D = $null;          E = "  +0.19%  " },
    @{ Row = 29; D = $null;          E = "  -0.21%  " },
    @{ Row = 30; D = "'1.25";        E = "  +6.83%  " },
    @{ Row = 32; D = $null;          E = "  +0.16%  " },
    @{ Row = 33; D = $null;          E = "  -2.17%  " },
    @{ Row = 34; D = "1.426.78";     E = "  -7.78%  " },
    @{ Row = 35; D = "'1.61";        E = "  -2.80%  " },
    @{ Row = 36; D = $null;          E = "  -0.64%  " },
    @{ Row = 37; D = $null;          E = "  +0.03%  " },
    @{ Row = 38; D = $null;          E = "  -2.16%  " },
    @{ Row = 39; D = "'0.0172";      E = "  +0.17%  " },
    @{ Row = 40; D = $null;          E = "  -2.88%  " },
    @{ Row = 41; D = "'69.38";       E = "  +0.23%  " },
    @{ Row = 42; D = $null;          E = "  -0.25%  " },
    @{ Row = 43; D = "'5.49";        E = "  -4.84%  " },
    @{ Row = 44; D = "1.833.28";     E = "  +0.15%  " },
    @{ Row = 45; D = $null;          E = "  -1.18%  " },
    @{ Row = 46; D = "'0.795";       E = "  +0.43%  " },
    @{ Row = 47; D = $null;          E = "  +5.48%  " },
    @{ Row = 48; D = "'89.22";       E = "  -0.88%  " },
    @{ Row = 49; D = "0.0₆0110";     E = "  -1.36%  " },
    @{ Row = 50; D = $null;          E = "  -1.06%  " },
    @{ Row = 51; D = "'7.91";        E = "  -3.76%  " }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
